# Insert a new data row at row 758 (shifts existing rows 758-799 down to 759-800)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(758).Insert()

# Write the new row's values. Column A holds date-like text (e.g. "2026/12/29"),
# so force text formatting before assigning to avoid Excel auto-converting it
# to a date serial number, then restore the default "Normal" style so the
# cell's style index matches its untouched neighbours.
$ws.Range("A758").NumberFormat = "@"
$ws.Range("A758").Value = "2026/02/02"
$ws.Range("A758").Style = "Normal"

$ws.Range("B758").Value = "月"
$ws.Range("C758").Value = 13
$ws.Range("D758").Value = 201
